$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")
$ws.Activate()

# Extend the formatting (wrap-text cell style used by every data row) from
# the last existing row down into the four new rows before filling values.
$ws.Range("A186:C186").Copy()
$ws.Range("A187:C190").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New translation rows appended after the existing last row (186).
$ws.Range("A187").Value = "cs"
$ws.Range("B187").Value = "lab.wizard.title"
$ws.Range("C187").Value = "Průvodci"

$ws.Range("A188").Value = "cs"
$ws.Range("B188").Value = "lab.wizard.subtitle"
$ws.Range("C188").Value = "Veškeré užitečné postupy v aplikaci jsou řešené pomocí průvodců; tady je najdete."

$ws.Range("A189").Value = "cs"
$ws.Range("B189").Value = "lab.wizard.build.title"
$ws.Range("C189").Value = "Průvodce novým buildem"

$ws.Range("A190").Value = "cs"
$ws.Range("B190").Value = "lab.wizard.build.subtitle"
$ws.Range("C190").Value = "Tento průvodce vám pomůže zaevidovat nový build."

# Match the saved view's selection (scrolled down to the newly added rows).
$ws.Range("B184").Select() | Out-Null
